$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: extend the trailing blank styled rows (516-521) by copying format from the still-blank row 506
$ws.Range("A506:E506").Copy()
$ws.Range("A516:E521").PasteSpecial(-4122)

# Step 2: fill in the 15 new event rows (501-515) using the format of row 500 as a template
$ws.Range("A500:E500").Copy()
$ws.Range("A501:E515").PasteSpecial(-4122)

# Row 501
$ws.Range("A501").Value2 = 45814
$ws.Range("B501").Value2 = 'BLITZ & DONNER'
$ws.Range("C501").Value2 = 'Projekt X'
$ws.Range("D501").Value2 = 'Bochum'
$link501 = 'https://www.instagram.com/reel/DKZw7aPC10N/?igsh=a3k4NGVldXNvNDg1'
$ws.Range("E501").Value2 = $link501
$ws.Hyperlinks.Add($ws.Range("E501"), $link501, "", "", $link501)
$ws.Range("E500").Copy()
$ws.Range("E501").PasteSpecial(-4122)
$len501 = $link501.Length
$c1_501 = $ws.Range("E501").Characters(1, $len501 - 1)
$c1_501.Font.Underline = 2
$c1_501.Font.Color = 65280
$c2_501 = $ws.Range("E501").Characters($len501, 1)
$c2_501.Font.Underline = 2
$c2_501.Font.Color = 65280

# Row 502
$ws.Range("A502").Value2 = 45822
$ws.Range("B502").Value2 = 'PUMP OPENING'
$ws.Range("C502").Value2 = 'SNRS'
$ws.Range("D502").Value2 = 'Dortmund'
$link502 = 'https://www.instagram.com/reel/DJcCofisk0b/?igsh=eHJ3dmo4anRmYmg5'
$ws.Range("E502").Value2 = $link502
$ws.Hyperlinks.Add($ws.Range("E502"), $link502, "", "", $link502)
$ws.Range("E500").Copy()
$ws.Range("E502").PasteSpecial(-4122)
$len502 = $link502.Length
$c1_502 = $ws.Range("E502").Characters(1, $len502 - 1)
$c1_502.Font.Underline = 2
$c1_502.Font.Color = 65280
$c2_502 = $ws.Range("E502").Characters($len502, 1)
$c2_502.Font.Underline = 2
$c2_502.Font.Color = 65280

# Row 503
$ws.Range("A503").Value2 = 45975
$ws.Range("B503").Value2 = 'PUMP CHURCH RAVE PRES. KUKO'
$ws.Range("C503").Value2 = 'Kreuzeskirche'
$ws.Range("D503").Value2 = 'Essen'
$link503 = 'https://www.instagram.com/reel/DIg2gvgMye0/?igsh=eXU2cDlwNzM4bXB4'
$ws.Range("E503").Value2 = $link503
$ws.Hyperlinks.Add($ws.Range("E503"), $link503, "", "", $link503)
$ws.Range("E500").Copy()
$ws.Range("E503").PasteSpecial(-4122)
$len503 = $link503.Length
$c1_503 = $ws.Range("E503").Characters(1, $len503 - 1)
$c1_503.Font.Underline = 2
$c1_503.Font.Color = 65280
$c2_503 = $ws.Range("E503").Characters($len503, 1)
$c2_503.Font.Underline = 2
$c2_503.Font.Color = 65280

# Row 504
$ws.Range("A504").Value2 = 45886
$ws.Range("B504").Value2 = 'PUMP ROOFTOP SESSION'
$ws.Range("C504").Value2 = 'Dortmunder U'
$ws.Range("D504").Value2 = 'Dortmund'
$link504 = 'https://www.instagram.com/reel/DJCWdFyM9mJ/?igsh=MTk1eDQxbnh4Y21icw=='
$ws.Range("E504").Value2 = $link504
$ws.Hyperlinks.Add($ws.Range("E504"), $link504, "", "", $link504)
$ws.Range("E500").Copy()
$ws.Range("E504").PasteSpecial(-4122)
$len504 = $link504.Length
$c1_504 = $ws.Range("E504").Characters(1, $len504 - 1)
$c1_504.Font.Underline = 2
$c1_504.Font.Color = 65280
$c2_504 = $ws.Range("E504").Characters($len504, 1)
$c2_504.Font.Underline = 2
$c2_504.Font.Color = 65280

# Row 505
$ws.Range("A505").Value2 = 45816
$ws.Range("B505").Value2 = 'PUMP'
$ws.Range("C505").Value2 = 'Oma Doris'
$ws.Range("D505").Value2 = 'Dortmund'
$link505 = 'https://www.instagram.com/reel/DJWS4z0MX35/?igsh=MTQxZTcwNTFubDJ0eA=='
# This link text already exists elsewhere in the sheet (E454); copy it directly so the shared string is reused
$ws.Range("E454").Copy()
$ws.Range("E505").PasteSpecial(-4163)
$ws.Range("E500").Copy()
$ws.Range("E505").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("E505"), $link505, "", "", $link505)
$ws.Range("EE454").Copy()
$ws.Range("E505").PasteSpecial(-4123)
$ws.Range("E500").Copy()
$ws.Range("E505").PasteSpecial(-4122)

# Row 506
$ws.Range("A506").Value2 = 45815
$ws.Range("B506").Value2 = 'SECRETRAVES'
$ws.Range("C506").Value2 = 'check event link'
$ws.Range("D506").Value2 = 'Köln'
$link506 = 'https://chat.whatsapp.com/B9bYDh2dFJ0CKGZzyH8gNK?fbclid=PAZXh0bgNhZW0CMTEAAafTo6-te2Yul4DhDMfLPF0BJdELrHZiMFLVX-7bSEngqpXwP4wK3MgisUplHA_aem_wwKjcPWJ9-CUbi75ZxsSHg&e=AT3VxQ3BI1gh4uA-t6K0K_FdxPt1K1x6tlppJfpMaeF-WgRDaBSRr9b8tzdbrlz480QX5ULYtAyXrEoFNrm3JpdklemJPH_HOdOfMMidqISlRjSZMveoo2nrkmHOmQ'
$ws.Range("E506").Value2 = $link506
$ws.Hyperlinks.Add($ws.Range("E506"), $link506, "", "", $link506)
$ws.Range("E500").Copy()
$ws.Range("E506").PasteSpecial(-4122)
$len506 = $link506.Length
$c1_506 = $ws.Range("E506").Characters(1, $len506 - 1)
$c1_506.Font.Underline = 2
$c1_506.Font.Color = 65280
$c2_506 = $ws.Range("E506").Characters($len506, 1)
$c2_506.Font.Underline = 2
$c2_506.Font.Color = 65280

# Row 507
$ws.Range("A507").Value2 = 45815
$ws.Range("B507").Value2 = 'BLACKCELL TECHNO BOILERROOM'
$ws.Range("C507").Value2 = 'Lessie Fair'
$ws.Range("D507").Value2 = 'Aachen'
$link507 = 'https://www.instagram.com/reel/DKRgoXaMpUm/?igsh=MXI0bHh4dXRxc3g4MA=='
$ws.Range("E507").Value2 = $link507
$ws.Hyperlinks.Add($ws.Range("E507"), $link507, "", "", $link507)
$ws.Range("E500").Copy()
$ws.Range("E507").PasteSpecial(-4122)
$len507 = $link507.Length
$c1_507 = $ws.Range("E507").Characters(1, $len507 - 1)
$c1_507.Font.Underline = 2
$c1_507.Font.Color = 65280
$c2_507 = $ws.Range("E507").Characters($len507, 1)
$c2_507.Font.Underline = 2
$c2_507.Font.Color = 65280

# Row 508
$ws.Range("A508").Value2 = 45821
$ws.Range("B508").Value2 = 'EHRENLOS RAVE PRES. FUNCHAINED'
$ws.Range("C508").Value2 = 'Stollen134'
$ws.Range("D508").Value2 = 'Dortmund'
$link508 = 'https://www.instagram.com/reel/DKST8YGoj_N/?igsh=NGg2aGNmNGNnb3Ju'
$ws.Range("E508").Value2 = $link508
$ws.Hyperlinks.Add($ws.Range("E508"), $link508, "", "", $link508)
$ws.Range("E500").Copy()
$ws.Range("E508").PasteSpecial(-4122)
$len508 = $link508.Length
$c1_508 = $ws.Range("E508").Characters(1, $len508 - 1)
$c1_508.Font.Underline = 2
$c1_508.Font.Color = 65280
$c2_508 = $ws.Range("E508").Characters($len508, 1)
$c2_508.Font.Underline = 2
$c2_508.Font.Color = 65280

# Row 509
$ws.Range("A509").Value2 = 45815
$ws.Range("B509").Value2 = 'BEHAVE THE RAVE'
$ws.Range("C509").Value2 = 'Franz Club'
$ws.Range("D509").Value2 = 'Paderborn'
$link509 = 'https://www.instagram.com/reel/DKb0OBTI80R/?igsh=NGhocDIxMXJoY254'
$ws.Range("E509").Value2 = $link509
$ws.Hyperlinks.Add($ws.Range("E509"), $link509, "", "", $link509)
$ws.Range("E500").Copy()
$ws.Range("E509").PasteSpecial(-4122)
$len509 = $link509.Length
$c1_509 = $ws.Range("E509").Characters(1, $len509 - 1)
$c1_509.Font.Underline = 2
$c1_509.Font.Color = 65280
$c2_509 = $ws.Range("E509").Characters($len509, 1)
$c2_509.Font.Underline = 2
$c2_509.Font.Color = 65280

# Row 510
$ws.Range("A510").Value2 = 45975
$ws.Range("B510").Value2 = 'CHROME'
$ws.Range("C510").Value2 = 'Bootshaus'
$ws.Range("D510").Value2 = 'Köln'
$link510 = 'https://www.instagram.com/chromecologne?igsh=MXQ3dzdjZTlyMTJ0Zw=='
$ws.Range("E510").Value2 = $link510
$ws.Hyperlinks.Add($ws.Range("E510"), $link510, "", "", $link510)
$ws.Range("E500").Copy()
$ws.Range("E510").PasteSpecial(-4122)
$len510 = $link510.Length
$c1_510 = $ws.Range("E510").Characters(1, $len510 - 1)
$c1_510.Font.Underline = 2
$c1_510.Font.Color = 65280
$c2_510 = $ws.Range("E510").Characters($len510, 1)
$c2_510.Font.Underline = 2
$c2_510.Font.Color = 65280

# Row 511
$ws.Range("A511").Value2 = 45836
$ws.Range("B511").Value2 = '24 STUNDEN RAVE'
$ws.Range("C511").Value2 = 'Essigfabrik & Elektroküche'
$ws.Range("D511").Value2 = 'Köln'
$link511 = 'https://www.instagram.com/reel/DKZhlp2sWGn/?igsh=MWFzeGlxZ25ubXlpcw=='
$ws.Range("E511").Value2 = $link511
$ws.Hyperlinks.Add($ws.Range("E511"), $link511, "", "", $link511)
$ws.Range("E500").Copy()
$ws.Range("E511").PasteSpecial(-4122)
$len511 = $link511.Length
$c1_511 = $ws.Range("E511").Characters(1, $len511 - 1)
$c1_511.Font.Underline = 2
$c1_511.Font.Color = 65280
$c2_511 = $ws.Range("E511").Characters($len511, 1)
$c2_511.Font.Underline = 2
$c2_511.Font.Color = 65280

# Row 512
$ws.Range("A512").Value2 = 45815
$ws.Range("B512").Value2 = 'TECHNO O PLOMO'
$ws.Range("C512").Value2 = 'Club 809'
$ws.Range("D512").Value2 = 'Essen'
$link512 = 'https://www.instagram.com/p/DJeq1KWMxId/?igsh=c3JoeHdkeXpwNW56'
$ws.Range("E512").Value2 = $link512
$ws.Hyperlinks.Add($ws.Range("E512"), $link512, "", "", $link512)
$ws.Range("E500").Copy()
$ws.Range("E512").PasteSpecial(-4122)
$len512 = $link512.Length
$c1_512 = $ws.Range("E512").Characters(1, $len512 - 1)
$c1_512.Font.Underline = 2
$c1_512.Font.Color = 65280
$c2_512 = $ws.Range("E512").Characters($len512, 1)
$c2_512.Font.Underline = 2
$c2_512.Font.Color = 65280

# Row 513
$ws.Range("A513").Value2 = 45822
$ws.Range("B513").Value2 = 'NIBIRII SOMMERFEST'
$ws.Range("C513").Value2 = 'Schrotty'
$ws.Range("D513").Value2 = 'Köln'
$link513 = 'https://www.instagram.com/p/DJW7jAQN57E/?igsh=aXJlbTdoeGc2NzVw'
$ws.Range("E513").Value2 = $link513
$ws.Hyperlinks.Add($ws.Range("E513"), $link513, "", "", $link513)
$ws.Range("E500").Copy()
$ws.Range("E513").PasteSpecial(-4122)
$len513 = $link513.Length
$c1_513 = $ws.Range("E513").Characters(1, $len513 - 1)
$c1_513.Font.Underline = 2
$c1_513.Font.Color = 65280
$c2_513 = $ws.Range("E513").Characters($len513, 1)
$c2_513.Font.Underline = 2
$c2_513.Font.Color = 65280

# Row 514
$ws.Range("A514").Value2 = 45849
$ws.Range("B514").Value2 = 'ANIMADO'
$ws.Range("C514").Value2 = 'Artheater'
$ws.Range("D514").Value2 = 'Köln'
$link514 = 'https://www.instagram.com/reel/DKZ_DkeIdtn/?igsh=ZWMwcXRqZ3BrNmVr'
$ws.Range("E514").Value2 = $link514
$ws.Hyperlinks.Add($ws.Range("E514"), $link514, "", "", $link514)
$ws.Range("E500").Copy()
$ws.Range("E514").PasteSpecial(-4122)
$len514 = $link514.Length
$c1_514 = $ws.Range("E514").Characters(1, $len514 - 1)
$c1_514.Font.Underline = 2
$c1_514.Font.Color = 65280
$c2_514 = $ws.Range("E514").Characters($len514, 1)
$c2_514.Font.Underline = 2
$c2_514.Font.Color = 65280

# Row 515
$ws.Range("A515").Value2 = 45836
$ws.Range("B515").Value2 = 'RAVE IM REINEKE'
$ws.Range("C515").Value2 = 'Reineke Fuchs'
$ws.Range("D515").Value2 = 'Köln'
$link515 = 'https://www.instagram.com/reel/DKbh4kBM_F3/?igsh=MTJyOTd4ejJiNGsycA=='
$ws.Range("E515").Value2 = $link515
$ws.Hyperlinks.Add($ws.Range("E515"), $link515, "", "", $link515)
$ws.Range("E500").Copy()
$ws.Range("E515").PasteSpecial(-4122)
$len515 = $link515.Length
$c1_515 = $ws.Range("E515").Characters(1, $len515 - 1)
$c1_515.Font.Underline = 2
$c1_515.Font.Color = 65280
$c2_515 = $ws.Range("E515").Characters($len515, 1)
$c2_515.Font.Underline = 2
$c2_515.Font.Color = 65280
